# "new tab for the justification of the initial states"
#
# Insert a new "Justification" column (right after the existing "Valeur EI"
# column, i.e. before the CT/LT justification & uncertainty columns) on both
# the "Niveau habitat" and "Niveau espece" sheets, then leave the UI focused
# on "Niveau espece" with the new column's second row selected.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Niveau habitat", "Niveau espece")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remember the width of column D so the freshly inserted column E can
    # be formatted the same way Excel formats a new column by default
    # (inheriting the look of the column immediately to its left).
    $leftWidth = $ws.Range("D1").ColumnWidth

    # Shift existing E:J one column to the right and create a new column E.
    $ws.Range("E1").EntireColumn.Insert()

    # Header for the new column.
    $ws.Range("E1").Value = "Justification"

    # Match column D's width on the new column.
    $ws.Range("E1").ColumnWidth = $leftWidth
}

# Activate "Niveau habitat" first, positioning its selection on the new
# column, then finish on "Niveau espece" (which becomes the active tab).
$wsHabitat = $wb.Worksheets.Item("Niveau habitat")
$wsHabitat.Activate()
[void]$wsHabitat.Range("E2").Select()

$wsEspece = $wb.Worksheets.Item("Niveau espece")
$wsEspece.Activate()
[void]$wsEspece.Range("E2").Select()
